# Edit: (1) swap the presentation's theme colour palette from the
# "Integral / Red Violet" scheme to the default "Office" scheme, and
# (2) re-style the three summary tables (slides 14-16) with the
# PowerPoint built-in table style {5B3F497D-C4FA-46F7-89C2-53C5C6FCCB56}
# instead of the custom local style {01E3AA2A-B34F-496A-963C-0BD6CF6E1852}.

$p = $ppt.ActivePresentation

# --- 1. Theme colours: Integral (Red Violet) -> Office -------------------
# ThemeColorScheme items are ordered dk1, lt1, dk2, lt2, accent1-6, hlink,
# folHlink. RGB is packed VBA-style (R + G*256 + B*65536), i.e. the reverse
# byte order of the hex "RRGGBB" string.
$cs = $p.SlideMaster.Theme.ThemeColorScheme
$cs.Item(1).RGB  = 0         # dk1      000000
$cs.Item(2).RGB  = 16777215  # lt1      FFFFFF
$cs.Item(3).RGB  = 6968388   # dk2      44546A
$cs.Item(4).RGB  = 15132391  # lt2      E7E6E6
$cs.Item(5).RGB  = 13998939  # accent1  5B9BD5
$cs.Item(6).RGB  = 3243501   # accent2  ED7D31
$cs.Item(7).RGB  = 10855845  # accent3  A5A5A5
$cs.Item(8).RGB  = 49407     # accent4  FFC000
$cs.Item(9).RGB  = 12874308  # accent5  4472C4
$cs.Item(10).RGB = 4697456   # accent6  70AD47
$cs.Item(11).RGB = 12673797  # hlink    0563C1
$cs.Item(12).RGB = 7491477   # folHlink 954F72

# --- 2. Table styles on the three comparison tables -----------------------
$newStyleId = "{5B3F497D-C4FA-46F7-89C2-53C5C6FCCB56}"
foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    foreach ($shp in $slide.Shapes) {
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle($newStyleId)
        }
    }
}
